# NextTransaction.xaml deleted and replaced with actions in transitions.
# The "Workblocks" config sheet had a Name/Value/Description row pair for
# the NextTransaction workblock (wbNextTransaction_Type /
# wbNextTransaction_SuppressSuccessful). Since that workblock no longer
# exists, remove its two rows entirely (shifting the remaining rows up).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Rows 11 and 12 hold the wbNextTransaction_Type / _SuppressSuccessful pair.
$ws.Rows("11:12").Delete()

# Restore the selection Excel leaves behind after such a row deletion.
$ws.Range("B21").Select()
